# Apply the "secondary scene with knight and flag" edit to the rubric.
#
# Semantics of the grading sheet (Sheet1):
#  - Column E on each feature row holds which milestone (I/II/III) the
#    feature was first completed on; column F holds "X" once it is
#    confirmed/checked off for credit.
#  - Row 83 is "Holy Grail Theme: everything must tie into the Monty
#    Python theme (knight, rabbit, etc...)" -- the author implemented a
#    secondary scene (knight + flag) for Milestone III, so that feature
#    is now marked complete on Milestone III.
#  - Row 91 ("Effective Use of GIT") is marked X as well.
#  - Several other features (rows 12, 29, 31, 68) that had been
#    (mistakenly) marked "III" were cleared out, since they were not
#    actually achieved.
# All of the totals in column G/H/I/J/K/L are formulas and recalculate
# automatically once the source cells change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out milestone markers that were mistakenly set.
$ws.Range("E12").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("E68").Value = ""

# Mark the "Holy Grail" secondary scene (knight + flag) as completed on
# Milestone III.
$ws.Range("E83").Value = "III"
$ws.Range("F83").Value = "X"

# Mark effective use of GIT as confirmed.
$ws.Range("E91").Value = "X"

# Update the cursor/viewport to reflect where the author was working.
$ws.Range("C1").Select() | Out-Null
$ws.Range("E12").Select() | Out-Null
